$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.565.24"
$ws.Range("E2").Value = "  +5.46%  "
$ws.Range("D3").Value = "1.724.22"
$ws.Range("E3").Value = "  +4.23%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'225.84"
$ws.Range("E5").Value = "  +3.28%  "
$ws.Range("D6").Value = "'0.5396"
$ws.Range("E6").Value = "  +3.05%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").Value = "'0.06613"
$ws.Range("E9").Value = "  +4.04%  "
$ws.Range("D10").Value = "'21.80"
$ws.Range("E10").Value = "  +6.35%  "
$ws.Range("D11").Value = "'0.07729"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "'4.618"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "1.730.83"
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("D14").Value = "1.962.61"
$ws.Range("E14").Value = "  +4.29%  "
$ws.Range("D15").Value = "'0.5865"
$ws.Range("E15").Value = "  +4.51%  "
$ws.Range("D16").Value = "0.0₅8314"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").Value = "'68.02"
$ws.Range("E17").Value = "  +3.82%  "
$ws.Range("D18").Value = "27.576.35"
$ws.Range("E18").Value = "  +5.54%  "
$ws.Range("D19").Value = "'221.08"
$ws.Range("E19").Value = "  +15.00%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "'4.736"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").Value = "'10.69"
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("D23").Value = "'6.099"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "'148.26"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("E26").Value = "  +3.44%  "
$ws.Range("E27").Value = "  +11.88%  "
$ws.Range("D28").Value = "'7.415"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("D29").Value = "'16.64"
$ws.Range("E29").Value = "  +4.27%  "
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("D31").Value = "'1.303"
$ws.Range("E31").Value = "  +2.58%  "
$ws.Range("E32").Value = "  +2.47%  "
$ws.Range("D33").Value = "'3.456"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").Value = "'1.662"
$ws.Range("E34").Value = "  +6.65%  "
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("D36").Value = "'2.825"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("D38").Value = "'0.5957"
$ws.Range("E38").Value = "  +5.21%  "
$ws.Range("E39").Value = "  +4.31%  "
$ws.Range("D40").Value = "'5.943"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8546"
$ws.Range("E41").Value = "  +2.39%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.056.54"
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "'101.51"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "1.867.75"
$ws.Range("E45").Value = "  +4.18%  "
$ws.Range("E46").Value = "  +4.92%  "
$ws.Range("E47").Value = "  +2.30%  "
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("D49").Value = "'0.4443"
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("E51").Value = "  +1.64%  "
